# Revert "Results 99% complete. Only conclusion left to do + appendices and
# random stuff": remove the two editorial sentences that were appended after
# "Expected Results Observed" and after the final "Passed" result.

$d = $word.ActiveDocument

# 1. Remove the "Despite this, the zipline..." commentary that followed
#    "Expected Results Observed" in the Test Record cell.
$zipline = ". Despite this, the zipline does not look very aesthetically pleasing at this point. The rope on the zipline does not connect properly, and the line is not straight, so the player phases through the line when getting closer to the bottom of the zipline. The zipline hook point could also do with a texture applied along with the ropes. However, functionally, the feature works fine. "
$found1 = $d.Content.Find.Execute($zipline, $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 2. Remove the ". More polish is required." commentary that followed the
#    final "Passed" result.
$polish = ". More polish is required."
$found2 = $d.Content.Find.Execute($polish, $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

Write-Output "zipline replaced: $found1"
Write-Output "polish replaced: $found2"
